$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ScenarioA")

# Row 8
$ws.Range("F8").Value = 17
$ws.Range("I8").Value = 8

# Row 9
$ws.Range("I9").Value = 8

# Row 10
$ws.Range("I10").Value = 8

# Row 11
$ws.Range("I11").Value = 8

# Row 12
$ws.Range("F12").Value = 22
$ws.Range("I12").Value = 8

# Row 13
$ws.Range("I13").Value = 8

# Row 14
$ws.Range("F14").Value = 17
$ws.Range("I14").Value = 8

# Row 15
$ws.Range("I15").Value = 8

# Row 16
$ws.Range("F16").Value = 80
$ws.Range("I16").Value = 8

# Row 17
$ws.Range("I17").Value = 8

# Row 18
$ws.Range("I18").Value = 8

# Update the active selection to match the diff (bottomRight pane -> H23)
$ws.Range("H23").Select()
